{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Applies the same textual changes described by the OOXML diff:\n//   - New page title / rewritten meta description\n//   - \"What we like\" bullet list rewritten (4 bullets)\n//   - \"What we don't like\" bullet list tweaks (2 bullets)\n\nconst replacements = [\n  {\n    from: \"Play Agent Destiny Slot for Free | Review of Play'n GO's Spy-Themed Game\",\n    to: \"Play Agent Destiny Free Slot Game\",\n  },\n  {\n    from: \"Fun spy theme inspired by classic British spies from the 60s\",\n    to: \"Spy theme straight out of the 60s\",\n  },\n  {\n    from: \"Features such as locking reels, colossal symbols & linked reels for more chances to win\",\n    to: \"Colossal symbols and locking reels for increased win potential\",\n  },\n  {\n    from: \"Free spin bonus with Mega Symbols and the Linked Reels feature for extra chances to win\",\n    to: \"Free spin bonus with potential for a maximum payout of 5000x\",\n  },\n  {\n    from: \"Wild symbol pays out 60x the winning bet when five matching symbols land on same payline\",\n    to: \"Specific graphic design inspired by classic British spies\",\n  },\n  {\n    from: \"Limited target audience of players who enjoy spy themed slots\",\n    to: \"Limited number of paylines (20)\",\n  },\n  {\n    from: \"No jackpot feature\",\n    to: \"No progressive jackpot feature\",\n  },\n  {\n    from: \"Read our review of Agent Destiny, an online slot game from Play'n GO. Play for free and enjoy features like locking reels and colossal symbols.\",\n    to: \"Read our review of Agent Destiny slot game and play for free.\",\n  },\n];\n\n// Use range-level search + replace (rather than whole-paragraph text\n// assignment) so untouched sibling runs in the same paragraph (e.g. the\n// leading empty `<w:r/>` Word keeps around some of these bullet items)\n// are left exactly as they were.\nfor (const { from, to } of replacements) {\n  const results = context.document.body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n#\n# Applies the same textual changes described by the OOXML diff:\n#   - New page title / rewritten meta description\n#   - \"What we like\" bullet list rewritten (4 bullets)\n#   - \"What we don't like\" bullet list tweaks (2 bullets)\n\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceAll = 2 (Word COM constants, spelled out\n# numerically since this host doesn't pre-seed the wd* enum names).\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll) | Out-Null\n}\n\nReplace-AllText \"Play Agent Destiny Slot for Free | Review of Play'n GO's Spy-Themed Game\" \"Play Agent Destiny Free Slot Game\"\nReplace-AllText \"Fun spy theme inspired by classic British spies from the 60s\" \"Spy theme straight out of the 60s\"\nReplace-AllText \"Features such as locking reels, colossal symbols & linked reels for more chances to win\" \"Colossal symbols and locking reels for increased win potential\"\nReplace-AllText \"Free spin bonus with Mega Symbols and the Linked Reels feature for extra chances to win\" \"Free spin bonus with potential for a maximum payout of 5000x\"\nReplace-AllText \"Wild symbol pays out 60x the winning bet when five matching symbols land on same payline\" \"Specific graphic design inspired by classic British spies\"\nReplace-AllText \"Limited target audience of players who enjoy spy themed slots\" \"Limited number of paylines (20)\"\nReplace-AllText \"No jackpot feature\" \"No progressive jackpot feature\"\nReplace-AllText \"Read our review of Agent Destiny, an online slot game from Play'n GO. Play for free and enjoy features like locking reels and colossal symbols.\" \"Read our review of Agent Destiny slot game and play for free.\"\n"}
